# Updated notebook, reran simulation
# Inserts two new fiber-orientation rows ("Holden" and "Rizzie Spiral") right after
# the "Spiral5" row, shifting all subsequent rows down by two, and renames the
# "Thomas Hex" entry to "Matthies Hex".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new blank rows at row 4 (pushes old rows 4..29 down to 6..31) ---
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(5).Insert()

# Fix up the formatting of the new index cells (column A) to match the existing
# bold/centered/bordered style used by every other index cell (copy format from row 3).
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# --- Populate new row 4: "Holden" ---
$ws.Range("A4").Value2 = 2
$ws.Range("B4").Value2 = "Holden"
$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W")
$row4vals = @(
    0.9986969927504069,
    0.9934181971878052,
    1.004350776807604,
    0.9934181971878052,
    1.001666314363113,
    1.00244529562576,
    0.9985869809728579,
    1.004350776807604,
    1.004350776807604,
    0.9957240877549126,
    1.001633859001183,
    1.004350776807604,
    1.001666314363113,
    0.997542255775459,
    1.000126647667985,
    0.9998117627861739,
    0.9978904975079254,
    0.9998117627861739,
    0.9995055673328449,
    1.000474609227797,
    0.9995653130579552
)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "4").Value2 = $row4vals[$i]
}

# --- Populate new row 5: "Rizzie Spiral" ---
$ws.Range("A5").Value2 = 3
$ws.Range("B5").Value2 = "Rizzie Spiral"
$row5vals = @(
    0.9939930582723508,
    0.9699374547468619,
    1.020026209035328,
    0.9699374547468619,
    1.007524345230188,
    1.011159863893649,
    0.9935408008312183,
    1.020026209035328,
    1.020026209035328,
    0.9804546667710523,
    1.007510796882044,
    1.020026209035328,
    1.007524345230188,
    0.9887308999885251,
    1.000532573030703,
    0.999162669670793,
    0.9903342002694228,
    0.999162669670793,
    0.9977572024608994,
    1.002211003775785,
    0.9980183994578368
)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "5").Value2 = $row5vals[$i]
}

# --- Rename "Thomas Hex" -> "Matthies Hex" ---
# After the insert, the row that used to hold "Thomas Hex" (old row 9) is now row 11.
$ws.Range("B11").Value2 = "Matthies Hex"

Write-Output "edit complete"
